$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 793-794, shifting existing rows 793-845 down to 795-847
$ws.Rows("793:794").Insert()

# Row 793
$ws.Cells.Item(793,1).Value = 6
$ws.Cells.Item(793,2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(793,3).Value = 'Metropolitana'
$ws.Cells.Item(793,4).Value = 44714
$ws.Cells.Item(793,5).Value = 13
$ws.Cells.Item(793,6).Value = 100112040
$ws.Cells.Item(793,7).Value = 'Cilantro'
$ws.Cells.Item(793,8).Value = 'Sin especificar'
$ws.Cells.Item(793,9).Value = 'Primera'
$ws.Cells.Item(793,10).Value = 620
$ws.Cells.Item(793,11).Value = 4000
$ws.Cells.Item(793,12).Value = 4500
$ws.Cells.Item(793,13).Value = 4218
$ws.Cells.Item(793,14).Value = '$/caja 36 atados'
$ws.Cells.Item(793,15).Value = 'Región Metropolitana'
$ws.Cells.Item(793,16).Value = 117
$ws.Cells.Item(793,17).Value = 36
$ws.Cells.Item(793,18).Value = 'Hortaliza'

# Row 794
$ws.Cells.Item(794,1).Value = 6
$ws.Cells.Item(794,2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(794,3).Value = 'Metropolitana'
$ws.Cells.Item(794,4).Value = 44714
$ws.Cells.Item(794,5).Value = 13
$ws.Cells.Item(794,6).Value = 100112040
$ws.Cells.Item(794,7).Value = 'Cilantro'
$ws.Cells.Item(794,8).Value = 'Sin especificar'
$ws.Cells.Item(794,9).Value = 'Primera'
$ws.Cells.Item(794,10).Value = 340
$ws.Cells.Item(794,11).Value = 7000
$ws.Cells.Item(794,12).Value = 8000
$ws.Cells.Item(794,13).Value = 7441
$ws.Cells.Item(794,14).Value = '$/docena de atados'
$ws.Cells.Item(794,15).Value = 'Región Metropolitana'
$ws.Cells.Item(794,16).Value = 2480
$ws.Cells.Item(794,17).Value = 3
$ws.Cells.Item(794,18).Value = 'Hortaliza'
